$d = $word.ActiveDocument

$pairs = @(
    @("957×5=4785", "302×5=1510"),
    @("752×8=6016", "823×2=1646"),
    @("559×2=1118", "663×6=3978"),
    @("114×4=456", "187×4=748"),
    @("845×6=5070", "377×2=754"),
    @("582×5=2910", "235×6=1410"),
    @("172×2=344", "331×9=2979"),
    @("967×7=6769", "515×3=1545"),
    @("195×2=390", "553×9=4977"),
    @("367×4=1468", "657×9=5913"),
    @("287×8=2296", "829×6=4974"),
    @("194×4=776", "520×4=2080"),
    @("939×9=8451", "504×6=3024"),
    @("363×7=2541", "516×9=4644"),
    @("324×9=2916", "780×8=6240"),
    @("136×2=272", "276×6=1656"),
    @("800×8=6400", "635×4=2540"),
    @("944×4=3776", "424×3=1272"),
    @("629×2=1258", "337×3=1011"),
    @("199×4=796", "758×7=5306"),
    @("259×6=1554", "401×4=1604"),
    @("506×2=1012", "704×5=3520"),
    @("335×2=670", "480×3=1440"),
    @("575×9=5175", "116×3=348"),
    @("653×7=4571", "603×7=4221")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
